Get-Variable | ForEach-Object { Write-Host $_.Name }
